$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("treatment")
$ws2 = $wb.Worksheets.Item("control")

# Sheet1 (treatment) updates
$ws1.Cells.Item(3, 2).Value = 0.731164928876849
$ws1.Cells.Item(3, 3).Value = 2.395201816281667
$ws1.Cells.Item(3, 4).Value = 5.903934177320373
$ws1.Cells.Item(4, 2).Value = 0.3082680437202807
$ws1.Cells.Item(4, 3).Value = 2.969666661165475
$ws1.Cells.Item(4, 4).Value = 8.574172475391292
$ws1.Cells.Item(5, 2).Value = 0.5921061238382785
$ws1.Cells.Item(5, 3).Value = 1.837761509389397
$ws1.Cells.Item(5, 4).Value = 3.122708225924273
$ws1.Cells.Item(7, 2).Value = 15.67984689656034
$ws1.Cells.Item(7, 3).Value = 41.46713634734983
$ws1.Cells.Item(7, 4).Value = 90.87829681474716
$ws1.Cells.Item(8, 2).Value = 0.3414679978567267
$ws1.Cells.Item(8, 3).Value = 2.24982902269246
$ws1.Cells.Item(8, 4).Value = 5.919526529873998
$ws1.Cells.Item(9, 2).Value = 0.6148592820564662
$ws1.Cells.Item(9, 3).Value = 1.578248869757236
$ws1.Cells.Item(9, 4).Value = 2.560025051636822
$ws1.Cells.Item(11, 2).Value = 0.4701008981129735
$ws1.Cells.Item(11, 3).Value = 1.717212382891513
$ws1.Cells.Item(11, 4).Value = 4.336772884958548
$ws1.Cells.Item(12, 2).Value = 0.216303147275205
$ws1.Cells.Item(12, 3).Value = 2.686844709663296
$ws1.Cells.Item(12, 4).Value = 8.18855655038591
$ws1.Cells.Item(13, 2).Value = 0.5004007210985406
$ws1.Cells.Item(13, 3).Value = 1.763631271272005
$ws1.Cells.Item(13, 4).Value = 3.078861719057596

$ws1.Cells.Item(14, 1).Value = "Speed meta analysis"
$ws1.Cells.Item(15, 1).Value = "mean (km/day)"
$ws1.Cells.Item(15, 2).Value = 3.653954736902968
$ws1.Cells.Item(15, 3).Value = 4.565123823551835
$ws1.Cells.Item(15, 4).Value = 5.632669521343581
$ws1.Cells.Item(16, 1).Value = "CoV² (RVAR)"
$ws1.Cells.Item(16, 2).Value = 0.02927619797602131
$ws1.Cells.Item(16, 3).Value = 0.1069301514646178
$ws1.Cells.Item(16, 4).Value = 0.2339562486043259
$ws1.Cells.Item(17, 1).Value = "CoV  (RSTD)"
$ws1.Cells.Item(17, 2).Value = 0.1763395690530721
$ws1.Cells.Item(17, 3).Value = 0.3370097987646193
$ws1.Cells.Item(17, 4).Value = 0.4984938008943804

# Sheet2 (control) updates
$ws2.Cells.Item(3, 2).Value = 1.061264292682796
$ws2.Cells.Item(3, 3).Value = 3.25548656999521
$ws2.Cells.Item(3, 4).Value = 7.394079564918341
$ws2.Cells.Item(4, 2).Value = 0.07571293480294974
$ws2.Cells.Item(4, 3).Value = 1.270882579976724
$ws2.Cells.Item(4, 4).Value = 4.101791450929976
$ws2.Cells.Item(5, 2).Value = 0.29912523165536
$ws2.Cells.Item(5, 3).Value = 1.225520896947137
$ws2.Cells.Item(5, 4).Value = 2.201683031906777
$ws2.Cells.Item(7, 2).Value = 0.3359549729327069
$ws2.Cells.Item(7, 3).Value = 0.6273037523207455
$ws2.Cells.Item(7, 4).Value = 1.075042756517226
$ws2.Cells.Item(8, 2).Value = 0.1135614307211665
$ws2.Cells.Item(8, 3).Value = 0.6434686353212727
$ws2.Cells.Item(8, 4).Value = 1.623707463961387
$ws2.Cells.Item(9, 2).Value = 0.352687610493711
$ws2.Cells.Item(9, 3).Value = 0.8395343445449107
$ws2.Cells.Item(9, 4).Value = 1.333609466123622
$ws2.Cells.Item(11, 2).Value = 0.2864639572642455
$ws2.Cells.Item(11, 3).Value = 0.6470016271740484
$ws2.Cells.Item(11, 4).Value = 1.234713535288357
$ws2.Cells.Item(12, 2).Value = 0.06738658327542106
$ws2.Cells.Item(12, 3).Value = 0.745958742022889
$ws2.Cells.Item(12, 4).Value = 2.219951378427903
$ws2.Cells.Item(13, 2).Value = 0.2781855479991011
$ws2.Cells.Item(13, 3).Value = 0.9255614329974937
$ws2.Cells.Item(13, 4).Value = 1.596686219911642

$ws2.Cells.Item(14, 1).Value = "Speed meta analysis"
$ws2.Cells.Item(15, 1).Value = "mean (km/day)"
$ws2.Cells.Item(15, 2).Value = 3.905356756047062
$ws2.Cells.Item(15, 3).Value = 4.876102998580794
$ws2.Cells.Item(15, 4).Value = 6.013107772340108
$ws2.Cells.Item(16, 1).Value = "CoV² (RVAR)"
$ws2.Cells.Item(16, 2).Value = 0.01320485772568828
$ws2.Cells.Item(16, 3).Value = 0.07107566359616169
$ws2.Cells.Item(16, 4).Value = 0.1766824793264817
$ws2.Cells.Item(17, 1).Value = "CoV  (RSTD)"
$ws2.Cells.Item(17, 2).Value = 0.1200468071097577
$ws2.Cells.Item(17, 3).Value = 0.2785121924202664
$ws2.Cells.Item(17, 4).Value = 0.4391172448707668
